$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 (logistic_embeddings)
$ws.Range("C5").Value = 0.312
$ws.Range("D5").Value = 0.441
$ws.Range("E5").Value = 0.48
$ws.Range("F5").Value = 0.54
$ws.Range("G5").Value = 0.5570000000000001
$ws.Range("H5").Value = 0.575

# Row 7 (classical-best-embeddings -> classical-best-embed)
$ws.Range("A7").Value = "classical-best-embed"
$ws.Range("C7").Value = 0.312
$ws.Range("E7").Value = 0.48
$ws.Range("F7").Value = 0.54

# Row 8 (BERT-base)
$ws.Range("C8").Value = 0.305
$ws.Range("D8").Value = 0.522
$ws.Range("F8").Value = 0.607
$ws.Range("G8").Value = 0.64
$ws.Range("H8").Value = 0.654

# Row 9 (BERT-base-nli)
$ws.Range("B9").Value = 0.266
$ws.Range("C9").Value = 0.412
$ws.Range("D9").Value = 0.547
$ws.Range("E9").Value = 0.575
$ws.Range("F9").Value = 0.591
$ws.Range("G9").Value = 0.624
$ws.Range("H9").Value = 0.638
